$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price values stored as text (inlineStr in source); use
# NumberFormat "@" + Style reset so Excel keeps them as text instead of
# silently converting to a number / leaving a quote-prefix style behind.
$textValues = @{
    "D2" = "242.58"
    "D4" = "5.398"
    "D6" = "3.408"
    "D7" = "6.277"
    "D8" = "1.125"
    "D9" = "0.8054"
    "D10" = "0.1420"
    "D11" = "0.07294"
    "D12" = "0.03077"
    "D13" = "0.03094"
    "D14" = "0.09356"
    "D15" = "3.902"
    "D16" = "0.001593"
    "D17" = "0.04803"
    "D18" = "0.0005807"
    "D19" = "0.006275"
    "D20" = "0.0009973"
    "D21" = "0.004063"
    "D22" = "0.0001499"
    "D23" = "3.732"
    "D24" = "2.154"
    "D27" = "0.0003996"
    "D40" = "0.03811"
    "D41" = "0.006654"
    "D42" = "0.1049"
    "D43" = "0.002638"
    "D45" = "0.00005606"
    "D47" = "0.3898"
    "D49" = "0.00002099"
    "D50" = "0.01009"
}
foreach ($addr in $textValues.Keys) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $textValues[$addr]
    $c.Style = "Normal"
}

# Plain text cells (coin names, links, volume/rank labels).
$ws.Range("E8").Value = "7FTXTokenFTTBestin24h"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONE"
